# 自动更新Excel文件 - 按"开始时间"与当前日期重新计算"剩余"天数
# today is advanced by one day (2026-01-14 -> 2026-01-15); rows whose
# remaining days would run out to 0 get "refilled" (start date reset to
# today, remaining reset back to the total days).

function Get-DayNum {
    param([int]$y, [int]$m, [int]$d)
    if ($m -le 2) {
        $y = $y - 1
    }
    $eraBase = $y
    if ($y -lt 0) { $eraBase = $y - 399 }
    $era = [Math]::Floor($eraBase / 400)
    $yoe = $y - $era * 400
    $mp = ($m + 9) % 12
    $doy = [Math]::Floor((153 * $mp + 2) / 5) + $d - 1
    $doe = $yoe * 365 + [Math]::Floor($yoe / 4) - [Math]::Floor($yoe / 100) + $doy
    return $era * 146097 + $doe - 719468
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = Get-DayNum 2026 1 15
$todaySerial = 20260115

$lastRow = 99
for ($r = 2; $r -le $lastRow; $r++) {
    $fCell = $ws.Cells.Item($r, 6)
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $fVal = $fCell.Value()
    $dVal = $dCell.Value()
    if ($fVal -eq $null -or $dVal -eq $null) {
        continue
    }

    $fStr = [string]([int64]$fVal)
    if ($fStr.Length -ne 8) {
        # malformed start date (e.g. row 36) - leave untouched
        continue
    }

    $fy = [int]$fStr.Substring(0, 4)
    $fm = [int]$fStr.Substring(4, 2)
    $fd = [int]$fStr.Substring(6, 2)

    $startDay = Get-DayNum $fy $fm $fd
    $totalDays = [int]$dVal
    $elapsed = $today - $startDay
    $remaining = $totalDays - $elapsed

    if ($remaining -le 0) {
        # ran out - refill: reset start date to today, remaining back to total
        $fCell.Value = $todaySerial
        $eCell.Value = $totalDays
    } else {
        $eCell.Value = $remaining
    }
}
